$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Structural edit: insert a blank row at 15 (shifts the vertical-timing block
#     down by one), and remove the obsolete "offset" row (was row 28, now row 29
#     after the insert). This matches how the block was re-laid out: a new blank
#     spacer row appears above "Vertical timing", and the old per-mode "offset"
#     row (style 7 / italic) is dropped entirely now that the Hx/Vx formulas
#     below compute directly from the porch/sync/back-porch values.
$ws.Rows.Item(15).Insert()
$ws.Rows.Item(29).Delete()

# --- Fix up the horizontal timing helper formulas (rows 25-28: Ha/Hb/Hc/Hd).
#     These used to chain off each other (and off the now-deleted "offset" row
#     B28); they are rewritten to compute straight from the porch/sync/back
#     porch inputs so Hsync/Hd no longer land on the wrong part of the scan.
$ws.Range("B25:E25").Formula = "=B12-1"
$ws.Range("B26:E26").Formula = "=B14-B11-B10-1"
$ws.Range("B27:E27").Formula = "=B14-B11-1"
$ws.Range("B28:E28").Formula = "=B14-1"

# --- Vertical timing helper formulas (rows 30-33: Va/Vb/Vc/Vd). Same fix as
#     above, now driven off the vertical porch/sync/back porch rows (18-20)
#     and the "Whole screen" row (22). The row labels also shift by one
#     (the old row 31/32/33 meanings slide down into 32/33, and a "Vc" label
#     now occupies row 32) without the row count itself changing.
$ws.Range("B30:E30").Formula = "=B20-1"

$ws.Range("A31").Value2 = "Vb"
$ws.Range("B31:E31").Formula = "=B22-B19-B18-1"

$ws.Range("A32").Value2 = "Vc"
$ws.Range("B32:E32").Formula = "=B22-B19-1"

$ws.Range("A33").Value2 = "Vd"
$ws.Range("B33:E33").Formula = "=B22-1"

# --- The old helper columns (G24:G27, G30:G33) and the stray explicit values
#     in F30:F33 are no longer needed now that the formulas above are
#     self-documenting; clear them out.
$ws.Range("G25:G28").ClearContents()
$ws.Range("F30:F33").ClearContents()
$ws.Range("G30:G33").ClearContents()

# --- Row label text for the "VGA controller settings" header moved down one
#     row as part of the insert above; nothing else to do there since Insert()
#     already carried the text/style with it.

# --- Update the selection to match the saved state of the edited workbook.
$ws.Range("H30").Select()
